$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 200
$ws.Cells.Item(200,1).Value = "WGE 195"
$ws.Cells.Item(200,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(200,3).Value = "13-02-2026"
$ws.Cells.Item(200,4).Value = 286962
$ws.Cells.Item(200,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(200,6).Value = 34413429360
$ws.Cells.Item(200,7).Value = "NEFT"
$ws.Cells.Item(200,8).Value = "SBIN0003229"
$ws.Cells.Item(200,9).Value = "AAAFW8862C"
$ws.Cells.Item(200,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(200,11).Value = "SAYAN BATTACHARYA"
$ws.Cells.Item(200,12).Value = "052cbf62-5587-47c0-9a01-cebe150b096e"
$ws.Cells.Item(200,13).Value = "ACC-14810110034736"
$ws.Cells.Item(200,14).Value = "UCBA0001481"
$ws.Cells.Item(200,21).Value = "pending"
$ws.Cells.Item(200,22).Value = 6500
$ws.Cells.Item(200,24).Value = "UT TEST (06 JAN 2026) RPA_ID : f792ad6849"
$ws.Cells.Item(200,25).Value = "HPCL, DUMKA"
$ws.Cells.Item(200,26).Value = "SITE EXPENSE"
$ws.Cells.Item(200,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(200,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(200,29).Value = 0
$ws.Cells.Item(200,30).Value = 0
$ws.Cells.Item(200,31).Value = 0

# Row 201
$ws.Cells.Item(201,1).Value = "WGE 195"
$ws.Cells.Item(201,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(201,3).Value = "13-02-2026"
$ws.Cells.Item(201,4).Value = 286962
$ws.Cells.Item(201,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(201,6).Value = 34413429360
$ws.Cells.Item(201,7).Value = "NEFT"
$ws.Cells.Item(201,8).Value = "SBIN0003229"
$ws.Cells.Item(201,9).Value = "AAAFW8862C"
$ws.Cells.Item(201,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(201,11).Value = "SAYAN BATTACHARYA"
$ws.Cells.Item(201,12).Value = "671d1244-c5c6-4d2f-9466-fcbfe7c1ffd9"
$ws.Cells.Item(201,13).Value = "ACC-14810110034736"
$ws.Cells.Item(201,14).Value = "UCBA0001481"
$ws.Cells.Item(201,21).Value = "pending"
$ws.Cells.Item(201,22).Value = 3500
$ws.Cells.Item(201,24).Value = "PATCH MATERIAL PURCHASE (06 JAN 2026) RPA_ID : 8dfa3e1799"
$ws.Cells.Item(201,25).Value = "HPCL, DUMKA"
$ws.Cells.Item(201,26).Value = "SITE EXPENSE"
$ws.Cells.Item(201,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(201,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(201,29).Value = 0
$ws.Cells.Item(201,30).Value = 0
$ws.Cells.Item(201,31).Value = 0

# Row 202
$ws.Cells.Item(202,1).Value = "WGE 195"
$ws.Cells.Item(202,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(202,3).Value = "13-02-2026"
$ws.Cells.Item(202,4).Value = 286962
$ws.Cells.Item(202,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(202,6).Value = 34413429360
$ws.Cells.Item(202,7).Value = "NEFT"
$ws.Cells.Item(202,8).Value = "SBIN0003229"
$ws.Cells.Item(202,9).Value = "AAAFW8862C"
$ws.Cells.Item(202,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(202,11).Value = "SAYAN BATTACHARYA"
$ws.Cells.Item(202,12).Value = "5fe335ed-66b5-475d-9625-8ffc27bf2a23"
$ws.Cells.Item(202,13).Value = "ACC-14810110034736"
$ws.Cells.Item(202,14).Value = "UCBA0001481"
$ws.Cells.Item(202,21).Value = "pending"
$ws.Cells.Item(202,22).Value = 5000
$ws.Cells.Item(202,24).Value = "HOLIDAY TESTING (06 JAN 2026) RPA_ID : 6240e20ac3"
$ws.Cells.Item(202,25).Value = "HPCL, DUMKA"
$ws.Cells.Item(202,26).Value = "SITE EXPENSE"
$ws.Cells.Item(202,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(202,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(202,29).Value = 0
$ws.Cells.Item(202,30).Value = 0
$ws.Cells.Item(202,31).Value = 0

# Row 203
$ws.Cells.Item(203,1).Value = "WGE 56"
$ws.Cells.Item(203,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(203,3).Value = "13-02-2026"
$ws.Cells.Item(203,4).Value = 286962
$ws.Cells.Item(203,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(203,6).Value = 34413429360
$ws.Cells.Item(203,7).Value = "NEFT"
$ws.Cells.Item(203,8).Value = "SBIN0003229"
$ws.Cells.Item(203,9).Value = "AAAFW8862C"
$ws.Cells.Item(203,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(203,11).Value = "Akhil"
$ws.Cells.Item(203,12).Value = "afcee0d6-126d-44f8-b247-b974b339597e"
$ws.Cells.Item(203,13).Value = "ACC-852110110008274"
$ws.Cells.Item(203,21).Value = "pending"
$ws.Cells.Item(203,22).Value = 400
$ws.Cells.Item(203,24).Value = "FUEL FOR GRASS CUTTING MACHINE (15.01.2026) RPA_ID : 2517b05830"
$ws.Cells.Item(203,25).Value = "IOCL FEROKE"
$ws.Cells.Item(203,26).Value = "FUEL EXPENSE"
$ws.Cells.Item(203,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(203,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(203,29).Value = 0
$ws.Cells.Item(203,30).Value = 0
$ws.Cells.Item(203,31).Value = 0

# Row 204
$ws.Cells.Item(204,1).Value = "WGP015"
$ws.Cells.Item(204,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(204,3).Value = "13-02-2026"
$ws.Cells.Item(204,4).Value = 286962
$ws.Cells.Item(204,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(204,6).Value = 34413429360
$ws.Cells.Item(204,7).Value = "NEFT"
$ws.Cells.Item(204,8).Value = "SBIN0003229"
$ws.Cells.Item(204,9).Value = "AAAFW8862C"
$ws.Cells.Item(204,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(204,11).Value = "NARAYAN TECHNOLOGY"
$ws.Cells.Item(204,12).Value = "530755f8-7f10-426b-ab5e-4fdf67009fbd"
$ws.Cells.Item(204,13).Value = "ACC-510101007107793"
$ws.Cells.Item(204,14).Value = "UBIN0903621"
$ws.Cells.Item(204,21).Value = "pending"
$ws.Cells.Item(204,22).Value = 98894
$ws.Cells.Item(204,24).Value = "Being payment for consumbale purchase RPA_ID : d7d46170df"
$ws.Cells.Item(204,25).Value = "ONGC-Electrical GOA"
$ws.Cells.Item(204,26).Value = "SITE EXPENSE"
$ws.Cells.Item(204,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(204,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(204,29).Value = 0
$ws.Cells.Item(204,30).Value = 0
$ws.Cells.Item(204,31).Value = 0

# Row 205
$ws.Cells.Item(205,1).Value = "WGP008"
$ws.Cells.Item(205,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(205,3).Value = "13-02-2026"
$ws.Cells.Item(205,4).Value = 286962
$ws.Cells.Item(205,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(205,6).Value = 34413429360
$ws.Cells.Item(205,7).Value = "DCR"
$ws.Cells.Item(205,8).Value = "SBIN0003229"
$ws.Cells.Item(205,9).Value = "AAAFW8862C"
$ws.Cells.Item(205,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(205,11).Value = "GAYATHRI ELECTRICALS"
$ws.Cells.Item(205,12).Value = "0682d10e-150c-4295-9d96-bf00b1cb68e6"
$ws.Cells.Item(205,13).Value = "ACC-39177475703"
$ws.Cells.Item(205,14).Value = "SBIN0000512"
$ws.Cells.Item(205,21).Value = "pending"
$ws.Cells.Item(205,22).Value = 281859.52
$ws.Cells.Item(205,24).Value = "Being payment for consumbale purchase RPA_ID : f5cb604fc7"
$ws.Cells.Item(205,25).Value = "ONGC-Electrical GOA"
$ws.Cells.Item(205,26).Value = "SITE EXPENSE"
$ws.Cells.Item(205,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(205,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(205,29).Value = 0
$ws.Cells.Item(205,30).Value = 0
$ws.Cells.Item(205,31).Value = 0

# Row 206
$ws.Cells.Item(206,1).Value = "WGE 131"
$ws.Cells.Item(206,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(206,3).Value = "13-02-2026"
$ws.Cells.Item(206,4).Value = 286962
$ws.Cells.Item(206,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(206,6).Value = 34413429360
$ws.Cells.Item(206,7).Value = "NEFT"
$ws.Cells.Item(206,8).Value = "SBIN0003229"
$ws.Cells.Item(206,9).Value = "AAAFW8862C"
$ws.Cells.Item(206,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(206,11).Value = "SAJIN SATHAR"
$ws.Cells.Item(206,12).Value = "89c7c91b-00df-436f-831e-2825fd4b54ab"
$ws.Cells.Item(206,13).Value = "ACC-5507101003171"
$ws.Cells.Item(206,14).Value = "CNRB0005507"
$ws.Cells.Item(206,21).Value = "pending"
$ws.Cells.Item(206,22).Value = 345
$ws.Cells.Item(206,24).Value = "Being payment for consumbale purchase RPA_ID : 61386e1dd1"
$ws.Cells.Item(206,25).Value = "ONGC-Electrical GOA"
$ws.Cells.Item(206,26).Value = "SITE EXPENSE"
$ws.Cells.Item(206,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(206,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(206,29).Value = 0
$ws.Cells.Item(206,30).Value = 0
$ws.Cells.Item(206,31).Value = 0

# Row 207
$ws.Cells.Item(207,1).Value = "WGP011"
$ws.Cells.Item(207,2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(207,3).Value = "13-02-2026"
$ws.Cells.Item(207,4).Value = 286962
$ws.Cells.Item(207,5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(207,6).Value = 34413429360
$ws.Cells.Item(207,7).Value = "NEFT"
$ws.Cells.Item(207,8).Value = "SBIN0003229"
$ws.Cells.Item(207,9).Value = "AAAFW8862C"
$ws.Cells.Item(207,10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(207,11).Value = "SHREE BALAJI ELECTRICAL"
$ws.Cells.Item(207,12).Value = "a11feec2-9735-456e-94d3-b8a3cb4bd241"
$ws.Cells.Item(207,13).Value = "ACC-125006695576"
$ws.Cells.Item(207,14).Value = "CNRB0017203"
$ws.Cells.Item(207,21).Value = "pending"
$ws.Cells.Item(207,22).Value = 277842
$ws.Cells.Item(207,24).Value = "Being payment for consumbale purchase RPA_ID : 1dfcedb8b6"
$ws.Cells.Item(207,25).Value = "ONGC-Electrical GOA"
$ws.Cells.Item(207,26).Value = "SITE EXPENSE"
$ws.Cells.Item(207,27).Value = "estimation@westernidc.com"
$ws.Cells.Item(207,28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(207,29).Value = 0
$ws.Cells.Item(207,30).Value = 0
$ws.Cells.Item(207,31).Value = 0

